$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 102.3975143333333
$ws.Range("H2").Value = 307.192543
$ws.Range("I2").Value = 0.822753865655704
$ws.Range("J2").Value = 0.8227538656557041
$ws.Range("M2").Value = 102.3975143333333
$ws.Range("N2").Value = 307.192543
$ws.Range("O2").Value = 0.822753865655704
$ws.Range("P2").Value = 0.8227538656557041
$ws.Range("Q2").Value = 10485.2509416452
$ws.Range("R2").Value = 94367.25847480685
$ws.Range("S2").Value = 0.6769239234514042
$ws.Range("T2").Value = 0.6769239234514044
$ws.Range("G3").Value = 102.3975143333333
$ws.Range("H3").Value = 307.192543
$ws.Range("I3").Value = 0.822753865655704
$ws.Range("J3").Value = 0.8227538656557041
$ws.Range("O3").Value = 0.1536069072592176
$ws.Range("P3").Value = 0.1536069072592176
$ws.Range("Q3").Value = 1957.580555029451
$ws.Range("R3").Value = 17618.22499526505
$ws.Range("S3").Value = 0.1263806767389385
$ws.Range("T3").Value = 0.1263806767389385
$ws.Range("G4").Value = 102.3975143333333
$ws.Range("H4").Value = 307.192543
$ws.Range("I4").Value = 0.822753865655704
$ws.Range("J4").Value = 0.8227538656557041
$ws.Range("M4").Value = 2.600356333333333
$ws.Range("N4").Value = 7.801069
$ws.Range("O4").Value = 0.02089360507685526
$ws.Range("P4").Value = 0.02089360507685526
$ws.Range("Q4").Value = 266.2700249142741
$ws.Range("R4").Value = 2396.430224228467
$ws.Range("S4").Value = 0.0171902943444663
$ws.Range("T4").Value = 0.01719029434446631
$ws.Range("G5").Value = 102.3975143333333
$ws.Range("H5").Value = 307.192543
$ws.Range("I5").Value = 0.822753865655704
$ws.Range("J5").Value = 0.8227538656557041
$ws.Range("M5").Value = 0.341712
$ws.Range("N5").Value = 1.025136
$ws.Range("O5").Value = 0.002745622008223115
$ws.Range("P5").Value = 0.002745622008223116
$ws.Range("Q5").Value = 34.990459417872
$ws.Range("R5").Value = 314.914134760848
$ws.Range("S5").Value = 0.002258971120894945
$ws.Range("T5").Value = 0.002258971120894946
$ws.Range("I6").Value = 0.1536069072592176
$ws.Range("J6").Value = 0.1536069072592176
$ws.Range("M6").Value = 102.3975143333333
$ws.Range("N6").Value = 307.192543
$ws.Range("O6").Value = 0.822753865655704
$ws.Range("P6").Value = 0.8227538656557041
$ws.Range("Q6").Value = 1957.580555029451
$ws.Range("R6").Value = 17618.22499526505
$ws.Range("S6").Value = 0.1263806767389385
$ws.Range("T6").Value = 0.1263806767389385
$ws.Range("I7").Value = 0.1536069072592176
$ws.Range("J7").Value = 0.1536069072592176
$ws.Range("O7").Value = 0.1536069072592176
$ws.Range("P7").Value = 0.1536069072592176
$ws.Range("S7").Value = 0.02359508195774188
$ws.Range("T7").Value = 0.02359508195774188
$ws.Range("I8").Value = 0.1536069072592176
$ws.Range("J8").Value = 0.1536069072592176
$ws.Range("M8").Value = 2.600356333333333
$ws.Range("N8").Value = 7.801069
$ws.Range("O8").Value = 0.02089360507685526
$ws.Range("P8").Value = 0.02089360507685526
$ws.Range("Q8").Value = 49.71221252217389
$ws.Range("R8").Value = 447.409912699565
$ws.Range("S8").Value = 0.003209402057351224
$ws.Range("T8").Value = 0.003209402057351224
$ws.Range("I9").Value = 0.1536069072592176
$ws.Range("J9").Value = 0.1536069072592176
$ws.Range("M9").Value = 0.341712
$ws.Range("N9").Value = 1.025136
$ws.Range("O9").Value = 0.002745622008223115
$ws.Range("P9").Value = 0.002745622008223116
$ws.Range("Q9").Value = 6.53266606104
$ws.Range("R9").Value = 58.79399454936
$ws.Range("S9").Value = 0.0004217465051859949
$ws.Range("T9").Value = 0.0004217465051859949
$ws.Range("G10").Value = 2.600356333333333
$ws.Range("H10").Value = 7.801069
$ws.Range("I10").Value = 0.02089360507685526
$ws.Range("J10").Value = 0.02089360507685526
$ws.Range("M10").Value = 102.3975143333333
$ws.Range("N10").Value = 307.192543
$ws.Range("O10").Value = 0.822753865655704
$ws.Range("P10").Value = 0.8227538656557041
$ws.Range("Q10").Value = 266.2700249142741
$ws.Range("R10").Value = 2396.430224228467
$ws.Range("S10").Value = 0.0171902943444663
$ws.Range("T10").Value = 0.01719029434446631
$ws.Range("G11").Value = 2.600356333333333
$ws.Range("H11").Value = 7.801069
$ws.Range("I11").Value = 0.02089360507685526
$ws.Range("J11").Value = 0.02089360507685526
$ws.Range("O11").Value = 0.1536069072592176
$ws.Range("P11").Value = 0.1536069072592176
$ws.Range("Q11").Value = 49.71221252217389
$ws.Range("R11").Value = 447.409912699565
$ws.Range("S11").Value = 0.003209402057351224
$ws.Range("T11").Value = 0.003209402057351224
$ws.Range("G12").Value = 2.600356333333333
$ws.Range("H12").Value = 7.801069
$ws.Range("I12").Value = 0.02089360507685526
$ws.Range("J12").Value = 0.02089360507685526
$ws.Range("M12").Value = 2.600356333333333
$ws.Range("N12").Value = 7.801069
$ws.Range("O12").Value = 0.02089360507685526
$ws.Range("P12").Value = 0.02089360507685526
$ws.Range("Q12").Value = 6.761853060306779
$ws.Range("R12").Value = 60.856677542761
$ws.Range("S12").Value = 0.0004365427331075918
$ws.Range("T12").Value = 0.0004365427331075919
$ws.Range("G13").Value = 2.600356333333333
$ws.Range("H13").Value = 7.801069
$ws.Range("I13").Value = 0.02089360507685526
$ws.Range("J13").Value = 0.02089360507685526
$ws.Range("M13").Value = 0.341712
$ws.Range("N13").Value = 1.025136
$ws.Range("O13").Value = 0.002745622008223115
$ws.Range("P13").Value = 0.002745622008223116
$ws.Range("Q13").Value = 0.8885729633760001
$ws.Range("R13").Value = 7.997156670384
$ws.Range("S13").Value = 0.00005736594193013601
$ws.Range("T13").Value = 0.00005736594193013603
$ws.Range("G14").Value = 0.341712
$ws.Range("H14").Value = 1.025136
$ws.Range("I14").Value = 0.002745622008223115
$ws.Range("J14").Value = 0.002745622008223116
$ws.Range("M14").Value = 102.3975143333333
$ws.Range("N14").Value = 307.192543
$ws.Range("O14").Value = 0.822753865655704
$ws.Range("P14").Value = 0.8227538656557041
$ws.Range("Q14").Value = 34.990459417872
$ws.Range("R14").Value = 314.914134760848
$ws.Range("S14").Value = 0.002258971120894945
$ws.Range("T14").Value = 0.002258971120894946
$ws.Range("G15").Value = 0.341712
$ws.Range("H15").Value = 1.025136
$ws.Range("I15").Value = 0.002745622008223115
$ws.Range("J15").Value = 0.002745622008223116
$ws.Range("O15").Value = 0.1536069072592176
$ws.Range("P15").Value = 0.1536069072592176
$ws.Range("Q15").Value = 6.53266606104
$ws.Range("R15").Value = 58.79399454936
$ws.Range("S15").Value = 0.0004217465051859949
$ws.Range("T15").Value = 0.0004217465051859949
$ws.Range("G16").Value = 0.341712
$ws.Range("H16").Value = 1.025136
$ws.Range("I16").Value = 0.002745622008223115
$ws.Range("J16").Value = 0.002745622008223116
$ws.Range("M16").Value = 2.600356333333333
$ws.Range("N16").Value = 7.801069
$ws.Range("O16").Value = 0.02089360507685526
$ws.Range("P16").Value = 0.02089360507685526
$ws.Range("Q16").Value = 0.8885729633760001
$ws.Range("R16").Value = 7.997156670384
$ws.Range("S16").Value = 0.00005736594193013601
$ws.Range("T16").Value = 0.00005736594193013603
$ws.Range("G17").Value = 0.341712
$ws.Range("H17").Value = 1.025136
$ws.Range("I17").Value = 0.002745622008223115
$ws.Range("J17").Value = 0.002745622008223116
$ws.Range("M17").Value = 0.341712
$ws.Range("N17").Value = 1.025136
$ws.Range("O17").Value = 0.002745622008223115
$ws.Range("P17").Value = 0.002745622008223116
$ws.Range("Q17").Value = 0.116767090944
$ws.Range("R17").Value = 1.050903818496
$ws.Range("S17").Value = 0.000007538440212039133
$ws.Range("T17").Value = 0.000007538440212039136
